# Edit process blank names: prefix the blank-filter/blank-sol sample
# names in column A with "EC1_" to match the EC1_ naming convention
# used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "EC1_blank-filter1"
$ws.Range("A31").Value = "EC1_blank-sol1"
$ws.Range("A32").Value = "EC1_blank-filter2"
$ws.Range("A33").Value = "EC1_blank-sol2"
$ws.Range("A34").Value = "EC1_blank-filter3"
$ws.Range("A35").Value = "EC1_blank-sol3"

# Leave selection on the last edited cell, matching the recorded view state.
$ws.Range("A35").Select()
